# "correção nos dados e inicio da analise PNAD 2009"
#
# The sheet had two stray "section header" rows that carried only an A-column
# label and no numeric data (row 5 "situação do domicílio" and row 8
# "grandes regiões e unidades da federação"). Remove them entirely so every
# remaining label row lines up with its data, shifting everything below up.
#
# Also fix the mislabeled B2 sub-header: it read "unnamed: 1_level_1" (a
# leftover pandas artifact) and should read "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dados")

# Delete the higher-numbered row first so the lower row index stays valid.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()

$ws.Cells.Item(2, 2).Value = "total"
